$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2 = 2
    3 = 2
    4 = 4
    5 = 1
    6 = 1
    7 = 1
    8 = 1
    9 = 3
    10 = 0
    11 = 1
    12 = 2
    13 = 1
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 1
    20 = 2
    21 = 1
    22 = 1
    23 = 0
    24 = 0
    25 = 2
    26 = 1
    27 = 2
    28 = 0
    29 = 2
    30 = 1
    31 = 1
    32 = 2
    33 = 0
    34 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
